$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.11601931561875
$ws.Range("C2").Value = 6.064933174463656
$ws.Range("D2").Value = 10.80832145204679
$ws.Range("F2").Value = 55.01270400377355
$ws.Range("G2").Value = 3.751717750793297
$ws.Range("K2").Value = 13.11887176344386
$ws.Range("L2").Value = 10.98213803574132
$ws.Range("B3").Value = 17.05524242793263
$ws.Range("C3").Value = 5.999734391021541
$ws.Range("D3").Value = 10.6628974521364
$ws.Range("F3").Value = 53.8429345540741
$ws.Range("G3").Value = 3.755425741005651
$ws.Range("K3").Value = 13.11909336979767
$ws.Range("L3").Value = 10.9651620454239
$ws.Range("B4").Value = 17.0245745823268
$ws.Range("C4").Value = 5.957822854517259
$ws.Range("D4").Value = 10.57152734512496
$ws.Range("F4").Value = 53.11361427749759
$ws.Range("G4").Value = 3.757817100409683
$ws.Range("K4").Value = 13.12465451522965
$ws.Range("L4").Value = 10.95706919139572
$ws.Range("B5").Value = 17.01375922874371
$ws.Range("C5").Value = 5.940269308914678
$ws.Range("D5").Value = 10.5337862859676
$ws.Range("F5").Value = 52.81392287966712
$ws.Range("G5").Value = 3.758820542308858
$ws.Range("K5").Value = 13.12828241051463
$ws.Range("L5").Value = 10.95435893048163
$ws.Range("B6").Value = 17.01206517908241
$ws.Range("C6").Value = 5.937325862495473
$ws.Range("D6").Value = 10.52748925925846
$ws.Range("F6").Value = 52.76401803809009
$ws.Range("G6").Value = 3.758988914615108
$ws.Range("K6").Value = 13.12896696301556
$ws.Range("L6").Value = 10.9539444246209
$ws.Range("B7").Value = 17.02442190100284
$ws.Range("C7").Value = 5.957588043689042
$ws.Range("D7").Value = 10.5710203875912
$ws.Range("F7").Value = 53.10958220451722
$ws.Range("G7").Value = 3.757830515834022
$ws.Range("K7").Value = 13.12469793326521
$ws.Range("L7").Value = 10.95703025871318
$ws.Range("B8").Value = 17.09369015579117
$ws.Range("C8").Value = 6.04284162657036
$ws.Range("D8").Value = 10.75862016963302
$ws.Range("F8").Value = 54.61184172923682
$ws.Range("G8").Value = 3.752972543089152
$ws.Range("K8").Value = 13.1178211000347
$ws.Range("L8").Value = 10.97580213057443
$ws.Range("B9").Value = 17.2817208683403
$ws.Range("C9").Value = 6.195192843573292
$ws.Range("D9").Value = 11.10933390395251
$ws.Range("F9").Value = 57.45741191065778
$ws.Range("G9").Value = 3.744350341478528
$ws.Range("K9").Value = 13.14743019466385
$ws.Range("L9").Value = 11.03101535327647
$ws.Range("B10").Value = 17.45075616909981
$ws.Range("C10").Value = 6.298198832757455
$ws.Range("D10").Value = 11.35567481097404
$ws.Range("F10").Value = 59.47115799376962
$ws.Range("G10").Value = 3.738559392610596
$ws.Range("K10").Value = 13.19543109901701
$ws.Range("L10").Value = 11.08265230271639
$ws.Range("B11").Value = 17.53411287956044
$ws.Range("C11").Value = 6.343145707301785
$ws.Range("D11").Value = 11.4651088664133
$ws.Range("F11").Value = 60.36753678036967
$ws.Range("G11").Value = 3.736041414325954
$ws.Range("K11").Value = 13.22293216577624
$ws.Range("L11").Value = 11.10850700305302
$ws.Range("B12").Value = 17.56658277687243
$ws.Range("C12").Value = 6.359893045981132
$ws.Range("D12").Value = 11.50615684829694
$ws.Range("F12").Value = 60.70390928029092
$ws.Range("G12").Value = 3.735104529393912
$ws.Range("K12").Value = 13.23415549806269
$ws.Range("L12").Value = 11.11863338123225
$ws.Range("B13").Value = 17.55954998876426
$ws.Range("C13").Value = 6.356298323789941
$ws.Range("D13").Value = 11.49733410356078
$ws.Range("F13").Value = 60.63160542205962
$ws.Range("G13").Value = 3.735305566992451
$ws.Range("K13").Value = 13.23170246066779
$ws.Range("L13").Value = 11.11643762418621
$ws.Range("B14").Value = 17.53676621804737
$ws.Range("C14").Value = 6.344528983568114
$ws.Range("D14").Value = 11.46849385201693
$ws.Range("F14").Value = 60.39527310395084
$ws.Range("G14").Value = 3.735964003817876
$ws.Range("K14").Value = 13.22383933286924
$ws.Range("L14").Value = 11.10933340316632
$ws.Range("B15").Value = 17.52292752936388
$ws.Range("C15").Value = 6.337284413424777
$ws.Range("D15").Value = 11.45077685431287
$ws.Range("F15").Value = 60.2501066319701
$ws.Range("G15").Value = 3.736369476302648
$ws.Range("K15").Value = 13.21912814761553
$ws.Range("L15").Value = 11.10502545440199
$ws.Range("B16").Value = 17.44543657628986
$ws.Range("C16").Value = 6.295223141775547
$ws.Range("D16").Value = 11.34846891296557
$ws.Range("F16").Value = 59.41216087779664
$ws.Range("G16").Value = 3.738726280102996
$ws.Range("K16").Value = 13.19374741141601
$ws.Range("L16").Value = 11.08100986609702
$ws.Range("B17").Value = 17.39953649073968
$ws.Range("C17").Value = 6.268931678146318
$ws.Range("D17").Value = 11.28502282163094
$ws.Range("F17").Value = 58.89289350989574
$ws.Range("G17").Value = 3.740201823533925
$ws.Range("K17").Value = 13.17962478028607
$ws.Range("L17").Value = 11.06687984653038
$ws.Range("B18").Value = 17.37374579449265
$ws.Range("C18").Value = 6.253629539782203
$ws.Range("D18").Value = 11.24828380936031
$ws.Range("F18").Value = 58.59238495523746
$ws.Range("G18").Value = 3.741061475543871
$ws.Range("K18").Value = 13.17203546710088
$ws.Range("L18").Value = 11.05897550330909
$ws.Range("B19").Value = 17.36511896264204
$ws.Range("C19").Value = 6.248417533498684
$ws.Range("D19").Value = 11.23580272505735
$ws.Range("F19").Value = 58.49032941221522
$ws.Range("G19").Value = 3.741354424504518
$ws.Range("K19").Value = 13.16955765129557
$ws.Range("L19").Value = 11.05633763028371
$ws.Range("B20").Value = 17.40435969147147
$ws.Range("C20").Value = 6.271749064407267
$ws.Range("D20").Value = 11.29180238467148
$ws.Range("F20").Value = 58.94836257137727
$ws.Range("G20").Value = 3.740043616075563
$ws.Range("K20").Value = 13.18107296206239
$ws.Range("L20").Value = 11.06836097471057
$ws.Range("B21").Value = 17.54343402190223
$ws.Range("C21").Value = 6.347993321089124
$ws.Range("D21").Value = 11.47697569280691
$ws.Range("F21").Value = 60.46477466458149
$ws.Range("G21").Value = 3.735770154828759
$ws.Range("K21").Value = 13.22612701052935
$ws.Range("L21").Value = 11.11141100886662
$ws.Range("B22").Value = 17.63958655388773
$ws.Range("C22").Value = 6.396232525845997
$ws.Range("D22").Value = 11.59570412543549
$ws.Range("F22").Value = 61.43786456890692
$ws.Range("G22").Value = 3.733074015816662
$ws.Range("K22").Value = 13.26028601610674
$ws.Range("L22").Value = 11.14150149557816
$ws.Range("B23").Value = 17.5877955592246
$ws.Range("C23").Value = 6.370631335153027
$ws.Range("D23").Value = 11.53255096415517
$ws.Range("F23").Value = 60.92022650231973
$ws.Range("G23").Value = 3.734504174537075
$ws.Range("K23").Value = 13.24162556220553
$ws.Range("L23").Value = 11.12526429982646
$ws.Range("B24").Value = 17.40217725832931
$ws.Range("C24").Value = 6.270475906130078
$ws.Range("D24").Value = 11.28873816233235
$ws.Range("F24").Value = 58.92329114174102
$ws.Range("G24").Value = 3.740115106329196
$ws.Range("K24").Value = 13.18041658793539
$ws.Range("L24").Value = 11.06769067360209
$ws.Range("B25").Value = 17.22535195029981
$ws.Range("C25").Value = 6.155553733393596
$ws.Range("D25").Value = 11.01638533605263
$ws.Range("F25").Value = 56.70003242054506
$ws.Range("G25").Value = 3.746586842017204
$ws.Range("K25").Value = 13.13480347076997
$ws.Range("L25").Value = 11.01412198967354
